$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

# Insert a new row 8 by duplicating the formatting of row 7 (the previous
# last data row), then overwrite the values with the new record so the
# cell styles (s="4" for A:E, s="5" for F:T) match the existing table.
$ws.Rows.Item($row).Insert()
$ws.Range("A7:T7").Copy()
$ws.Range("A8:T8").PasteSpecial(-4104)  # xlPasteAll

$ws.Cells.Item($row, 1).Value = 92634000
$ws.Cells.Item($row, 2).Value = "Лаишевский МР"
$ws.Cells.Item($row, 3).Value = 2020
$ws.Cells.Item($row, 4).Value = 3578
$ws.Cells.Item($row, 5).Value = 47423
$ws.Cells.Item($row, 6).Value = 0.34242877928431353
$ws.Cells.Item($row, 7).Value = 37980
$ws.Cells.Item($row, 8).Value = 0.54330177340109231
$ws.Cells.Item($row, 9).Value = 0.023722666216814625
$ws.Cells.Item($row, 10).Value = 53.024439617906921
$ws.Cells.Item($row, 11).Value = 38.1
$ws.Cells.Item($row, 12).Value = 0.0036901925226156085
$ws.Cells.Item($row, 13).Value = 0.0011176011639921558
$ws.Cells.Item($row, 14).Value = 0.01201526685363642
$ws.Cells.Item($row, 15).Value = 0.35423739535668347
$ws.Cells.Item($row, 16).Value = 2.7574299390591066
$ws.Cells.Item($row, 17).Value = 140.13062016321194
$ws.Cells.Item($row, 18).Value = 0.00088564620542774602
$ws.Cells.Item($row, 19).Value = 0.056723530776205636
$ws.Cells.Item($row, 20).Value = 873.40777892583776

$ws.Range("C13").Select()
